$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping the cell a plain text
# cell (matching the sheet's existing inline-string / text convention for
# the Price column), avoiding any float round-trip / style residue.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Price (column D) updates
Set-TextValue "D2"  "245.81"
Set-TextValue "D3"  "24.24"
Set-TextValue "D4"  "5.376"
Set-TextValue "D5"  "0.05732"
Set-TextValue "D6"  "6.510"
Set-TextValue "D7"  "3.142"
Set-TextValue "D8"  "0.8165"
Set-TextValue "D9"  "0.8688"
Set-TextValue "D11" "0.06985"
Set-TextValue "D14" "0.09368"
Set-TextValue "D15" "3.741"
Set-TextValue "D16" "0.001541"
Set-TextValue "D18" "0.0005998"
Set-TextValue "D19" "0.006165"
Set-TextValue "D20" "0.001242"
Set-TextValue "D21" "0.004781"
Set-TextValue "D23" "3.530"
Set-TextValue "D40" "0.03695"
Set-TextValue "D41" "0.006439"
Set-TextValue "D42" "0.1054"
Set-TextValue "D44" "0.008655"
Set-TextValue "D48" "0.002508"
Set-TextValue "D49" "0.00002099"
Set-TextValue "D50" "0.0001999"

# Volume(1h) (column E) updates - "Bestin24h" / "Worstin24h" labels moved to new rows
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
